# Bump the cached "datetimeFigureOut" auto-date field text from
# 11/30/18 -> 12/1/18 on every master/layout that carries one
# (Slide Master, Notes Master, Handout Master, and the two slide
# layouts that override the date placeholder: "Blank" and
# "1_Title Slide").

$p = $ppt.ActivePresentation
$newDate = "12/1/18"

# Slide Master - "Date Placeholder 3"
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# Notes Master - "Date Placeholder 2"
$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Handout Master - "Date Placeholder 2"
$handoutMaster = $p.HandoutMaster
$handoutMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Slide Layout 6 ("Blank") - "Date Placeholder 1"
$layout6 = $master.CustomLayouts.Item(6)
$layout6.Shapes.Item(1).TextFrame.TextRange.Text = $newDate

# Slide Layout 7 ("1_Title Slide") - "Date Placeholder 3"
$layout7 = $master.CustomLayouts.Item(7)
$layout7.Shapes.Item(3).TextFrame.TextRange.Text = $newDate
